$wb = $excel.ActiveWorkbook

# --- Statistics sheet updates (columns D, E, M, N for rows 2-15) ---
$statsWs = $wb.Worksheets.Item("Statistics")
$statsWs.Range("D2").Value = 22.575200000000002
$statsWs.Range("E2").Value = -0.3149233375000122
$statsWs.Range("M2").Value = 250.11335994715682
$statsWs.Range("N2").Value = 0.26231554296902004
$statsWs.Range("D3").Value = 20.640013500000016
$statsWs.Range("E3").Value = 8.683079560000024
$statsWs.Range("M3").Value = 84.45377269046827
$statsWs.Range("N3").Value = 0.9674177198595402
$statsWs.Range("D4").Value = 9.746518250000001
$statsWs.Range("E4").Value = -0.002256849999991317
$statsWs.Range("M4").Value = 113.13687585280597
$statsWs.Range("N4").Value = 2.6445662227193196
$statsWs.Range("D5").Value = 6.911729999999988
$statsWs.Range("E5").Value = -5.612780425000013
$statsWs.Range("M5").Value = 73.40306221563604
$statsWs.Range("N5").Value = -0.06948994440095646
$statsWs.Range("D6").Value = 48.13514399999998
$statsWs.Range("E6").Value = 2.6775433549999974
$statsWs.Range("M6").Value = 34.27482010611674
$statsWs.Range("N6").Value = 4.68549706847654
$statsWs.Range("D7").Value = 5.065843999999999
$statsWs.Range("E7").Value = 0.12528425499999862
$statsWs.Range("M7").Value = 44.198782873567325
$statsWs.Range("N7").Value = 0.5729206984451949
$statsWs.Range("D8").Value = 12.885081249999997
$statsWs.Range("E8").Value = -0.7474411775000007
$statsWs.Range("M8").Value = 158.95269213622393
$statsWs.Range("N8").Value = -0.43895271004780057
$statsWs.Range("D9").Value = 5.461013499999999
$statsWs.Range("E9").Value = -0.41697330750000017
$statsWs.Range("M9").Value = 88.95137547579543
$statsWs.Range("N9").Value = 0.7593701311933501
$statsWs.Range("D10").Value = 41.630899
$statsWs.Range("E10").Value = -0.14258761999998626
$statsWs.Range("M10").Value = 74.7849416793606
$statsWs.Range("N10").Value = 0.06169703426013484
$statsWs.Range("D11").Value = 15.39517
$statsWs.Range("E11").Value = -1.5722586850000013
$statsWs.Range("M11").Value = 163.61691374741625
$statsWs.Range("N11").Value = -0.15515968164643823
$statsWs.Range("D12").Value = 4.101573000000001
$statsWs.Range("E12").Value = -0.4640823424999958
$statsWs.Range("M12").Value = 26.104832130023166
$statsWs.Range("N12").Value = -0.12642862324927506
$statsWs.Range("D13").Value = 10.770157
$statsWs.Range("E13").Value = -1.4136232800000013
$statsWs.Range("M13").Value = 155.04844673000753
$statsWs.Range("N13").Value = 0.8101772820826909
$statsWs.Range("D14").Value = 18.460378249999998
$statsWs.Range("E14").Value = -1.4046855824999938
$statsWs.Range("M14").Value = 160.64476524940255
$statsWs.Range("N14").Value = -0.8240053563438892
$statsWs.Range("D15").Value = 221.77872175
$statsWs.Range("E15").Value = -0.6057054375000632
$statsWs.Range("M15").Value = 1427.684640833981
$statsWs.Range("N15").Value = 9.14992538431784

# --- Speeds sheet updates (columns B, C, H, I for rows 2-15) ---
$speedsWs = $wb.Worksheets.Item("Speeds")
$speedsWs.Range("B2").Value = 11.079120448419363
$speedsWs.Range("C2").Value = 0.08833257681409007
$speedsWs.Range("H2").Value = 39.88483361430971
$speedsWs.Range("I2").Value = 0.31799727653072424
$speedsWs.Range("B3").Value = 4.091749876542871
$speedsWs.Range("C3").Value = 0.10045046839319675
$speedsWs.Range("H3").Value = 14.730299555554337
$speedsWs.Range("I3").Value = 0.36162168621550833
$speedsWs.Range("B4").Value = 11.607927359373278
$speedsWs.Range("C4").Value = 0.14554312498087568
$speedsWs.Range("H4").Value = 41.7885384937438
$speedsWs.Range("I4").Value = 0.5239552499311525
$speedsWs.Range("B5").Value = 10.620070838362635
$speedsWs.Range("C5").Value = 1.3974906086895729
$speedsWs.Range("H5").Value = 38.232255018105484
$speedsWs.Range("I5").Value = 5.030966191282462
$speedsWs.Range("B6").Value = 0.7120539642743513
$speedsWs.Range("C6").Value = 0.023556205107802024
$speedsWs.Range("H6").Value = 2.563394271387665
$speedsWs.Range("I6").Value = 0.08480233838808729
$speedsWs.Range("B7").Value = 8.724860630048484
$speedsWs.Range("C7").Value = 0.06347725813862823
$speedsWs.Range("H7").Value = 31.409498268174545
$speedsWs.Range("I7").Value = 0.22851812929906162
$speedsWs.Range("B8").Value = 12.336180816572186
$speedsWs.Range("C8").Value = 0.1361980698348556
$speedsWs.Range("H8").Value = 44.410250939659875
$speedsWs.Range("I8").Value = 0.4903130514054802
$speedsWs.Range("B9").Value = 16.288437205986664
$speedsWs.Range("C9").Value = 0.24533003970386885
$speedsWs.Range("H9").Value = 58.63837394155199
$speedsWs.Range("I9").Value = 0.8831881429339279
$speedsWs.Range("B10").Value = 1.7963806565733929
$speedsWs.Range("C10").Value = 0.019302347978012876
$speedsWs.Range("H10").Value = 6.466970363664214
$speedsWs.Range("I10").Value = 0.06948845272084636
$speedsWs.Range("B11").Value = 10.627808185776205
$speedsWs.Range("C11").Value = 0.17788515947154618
$speedsWs.Range("H11").Value = 38.26010946879434
$speedsWs.Range("I11").Value = 0.6403865740975663
$speedsWs.Range("B12").Value = 6.3645903973970865
$speedsWs.Range("C12").Value = 0.3308513311523656
$speedsWs.Range("H12").Value = 22.91252543062951
$speedsWs.Range("I12").Value = 1.1910647921485162
$speedsWs.Range("B13").Value = 14.396117598843501
$speedsWs.Range("C13").Value = 0.2768486125091275
$speedsWs.Range("H13").Value = 51.826023355836604
$speedsWs.Range("I13").Value = 0.9966550050328591
$speedsWs.Range("B14").Value = 8.70213833508003
$speedsWs.Range("C14").Value = 0.11094794078322623
$speedsWs.Range("H14").Value = 31.327698006288106
$speedsWs.Range("I14").Value = 0.39941258681961445
$speedsWs.Range("B15").Value = 6.437428395152074
$speedsWs.Range("C15").Value = 0.05091387077458165
$speedsWs.Range("H15").Value = 23.17474222254747
$speedsWs.Range("I15").Value = 0.18328993478849395
